$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A81").Value = "Kamin"
$ws.Range("B81").Value = "Ziegel"

$ws.Range("A80:B80").Copy() | Out-Null
$ws.Range("A81:B81").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

foreach ($fc in $ws.Cells.FormatConditions) {
    $addr = $fc.AppliesTo.Address()
    if ($addr -eq "`$A`$2:`$A`$80") {
        $fc.ModifyAppliesToRange($ws.Range("A2:A81")) | Out-Null
    } elseif ($addr -eq "`$A`$2:`$B`$80") {
        $fc.ModifyAppliesToRange($ws.Range("A2:B81")) | Out-Null
    } elseif ($addr -eq "`$B`$2:`$B`$80") {
        $fc.ModifyAppliesToRange($ws.Range("B2:B81")) | Out-Null
    }
}

$ws.Range("B81").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 70

